$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,14
$data[0,0] = 27.67634766666667
$data[0,1] = 83.029043
$data[0,2] = 0.005965811625935536
$data[0,3] = 0.005965811625935536
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 0.055614
$data[0,7] = 0.166842
$data[0,8] = 0.003173711121411028
$data[0,9] = 0.003173711121411028
$data[0,10] = 1.539192399134
$data[0,11] = 13.852731592206
$data[0,12] = 0.00001893376270547482
$data[0,13] = 0.00001893376270547482
$data[1,0] = 27.67634766666667
$data[1,1] = 83.029043
$data[1,2] = 0.005965811625935536
$data[1,3] = 0.005965811625935536
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 0.2521723333333333
$data[1,7] = 0.756517
$data[1,8] = 0.01439065952479895
$data[1,9] = 0.01439065952479895
$data[1,10] = 6.979209169247889
$data[1,11] = 62.812882523231
$data[1,12] = 0.00008585196389792555
$data[1,13] = 0.00008585196389792555
$data[2,0] = 27.67634766666667
$data[2,1] = 83.029043
$data[2,2] = 0.005965811625935536
$data[2,3] = 0.005965811625935536
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 4.376294333333333
$data[2,7] = 13.128883
$data[2,8] = 0.2497409644382361
$data[2,9] = 0.2497409644382361
$data[2,10] = 121.1198434609965
$data[2,11] = 1090.078591148969
$data[2,12] = 0.001489907549117982
$data[2,13] = 0.001489907549117982
$data[3,0] = 27.67634766666667
$data[3,1] = 83.029043
$data[3,2] = 0.005965811625935536
$data[3,3] = 0.005965811625935536
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 12.83925333333333
$data[3,7] = 38.51776
$data[3,8] = 0.732694664915554
$data[3,9] = 0.7326946649155539
$data[3,10] = 355.3436390337422
$data[3,11] = 3198.09275130368
$data[3,12] = 0.004371118350214154
$data[3,13] = 0.004371118350214153
$data[4,0] = 42.300692
$data[4,1] = 126.902076
$data[4,2] = 0.009118181457976757
$data[4,3] = 0.009118181457976757
$data[4,4] = 2
$data[4,5] = 0.6666666666666666
$data[4,6] = 0.055614
$data[4,7] = 0.166842
$data[4,8] = 0.003173711121411028
$data[4,9] = 0.003173711121411028
$data[4,10] = 2.352510684888
$data[4,11] = 21.172596163992
$data[4,12] = 0.00002893847390022466
$data[4,13] = 0.00002893847390022465
$data[5,0] = 42.300692
$data[5,1] = 126.902076
$data[5,2] = 0.009118181457976757
$data[5,3] = 0.009118181457976757
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 0.2521723333333333
$data[5,7] = 0.756517
$data[5,8] = 0.01439065952479895
$data[5,9] = 0.01439065952479895
$data[5,10] = 10.66706420325467
$data[5,11] = 96.003577829292
$data[5,12] = 0.0001312166448470784
$data[5,13] = 0.0001312166448470784
$data[6,0] = 42.300692
$data[6,1] = 126.902076
$data[6,2] = 0.009118181457976757
$data[6,3] = 0.009118181457976757
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 4.376294333333333
$data[6,7] = 13.128883
$data[6,8] = 0.2497409644382361
$data[6,9] = 0.2497409644382361
$data[6,10] = 185.1202786956786
$data[6,11] = 1666.082508261108
$data[6,12] = 0.002277183431237957
$data[6,13] = 0.002277183431237957
$data[7,0] = 42.300692
$data[7,1] = 126.902076
$data[7,2] = 0.009118181457976757
$data[7,3] = 0.009118181457976757
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 12.83925333333333
$data[7,7] = 38.51776
$data[7,8] = 0.732694664915554
$data[7,9] = 0.7326946649155539
$data[7,10] = 543.1093007633067
$data[7,11] = 4887.98370686976
$data[7,12] = 0.006680842907991497
$data[7,13] = 0.006680842907991496
$data[8,0] = 29.593002
$data[8,1] = 88.779006
$data[8,2] = 0.006378958578792732
$data[8,3] = 0.006378958578792732
$data[8,4] = 2
$data[8,5] = 0.6666666666666666
$data[8,6] = 0.055614
$data[8,7] = 0.166842
$data[8,8] = 0.003173711121411028
$data[8,9] = 0.003173711121411028
$data[8,10] = 1.645785213228
$data[8,11] = 14.812066919052
$data[8,12] = 0.00002024497178453478
$data[8,13] = 0.00002024497178453478
$data[9,0] = 29.593002
$data[9,1] = 88.779006
$data[9,2] = 0.006378958578792732
$data[9,3] = 0.006378958578792732
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 0.2521723333333333
$data[9,7] = 0.756517
$data[9,8] = 0.01439065952479895
$data[9,9] = 0.01439065952479895
$data[9,10] = 7.462536364678
$data[9,11] = 67.162827282102
$data[9,12] = 0.00009179742103020163
$data[9,13] = 0.00009179742103020162
$data[10,0] = 29.593002
$data[10,1] = 88.779006
$data[10,2] = 0.006378958578792732
$data[10,3] = 0.006378958578792732
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 4.376294333333333
$data[10,7] = 13.128883
$data[10,8] = 0.2497409644382361
$data[10,9] = 0.2497409644382361
$data[10,10] = 129.507686958922
$data[10,11] = 1165.569182630298
$data[10,12] = 0.001593087267579257
$data[10,13] = 0.001593087267579257
$data[11,0] = 29.593002
$data[11,1] = 88.779006
$data[11,2] = 0.006378958578792732
$data[11,3] = 0.006378958578792732
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 12.83925333333333
$data[11,7] = 38.51776
$data[11,8] = 0.732694664915554
$data[11,9] = 0.7326946649155539
$data[11,10] = 379.95204957184
$data[11,11] = 3419.56844614656
$data[11,12] = 0.00467382891839874
$data[11,13] = 0.004673828918398739
$data[12,0] = 4539.588785666667
$data[12,1] = 13618.766357
$data[12,2] = 0.9785370483372949
$data[12,3] = 0.978537048337295
$data[12,4] = 2
$data[12,5] = 0.6666666666666666
$data[12,6] = 0.055614
$data[12,7] = 0.166842
$data[12,8] = 0.003173711121411028
$data[12,9] = 0.003173711121411028
$data[12,10] = 252.464690726066
$data[12,11] = 2272.182216534594
$data[12,12] = 0.003105593913020794
$data[12,13] = 0.003105593913020794
$data[13,0] = 4539.588785666667
$data[13,1] = 13618.766357
$data[13,2] = 0.9785370483372949
$data[13,3] = 0.978537048337295
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 0.2521723333333333
$data[13,7] = 0.756517
$data[13,8] = 0.01439065952479895
$data[13,9] = 0.01439065952479895
$data[13,10] = 1144.758696455397
$data[13,11] = 10302.82826809857
$data[13,12] = 0.01408179349502375
$data[13,13] = 0.01408179349502375
$data[14,0] = 4539.588785666667
$data[14,1] = 13618.766357
$data[14,2] = 0.9785370483372949
$data[14,3] = 0.978537048337295
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 4.376294333333333
$data[14,7] = 13.128883
$data[14,8] = 0.2497409644382361
$data[14,9] = 0.2497409644382361
$data[14,10] = 19866.57667837658
$data[14,11] = 178799.1901053893
$data[14,12] = 0.2443807861903009
$data[14,13] = 0.2443807861903009
$data[15,0] = 4539.588785666667
$data[15,1] = 13618.766357
$data[15,2] = 0.9785370483372949
$data[15,3] = 0.978537048337295
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 12.83925333333333
$data[15,7] = 38.51776
$data[15,8] = 0.732694664915554
$data[15,9] = 0.7326946649155539
$data[15,10] = 58284.93044833337
$data[15,11] = 524564.3740350003
$data[15,12] = 0.7169688747389495
$data[15,13] = 0.7169688747389495

$range = $ws.Range("G2:T17")
$range.Value2 = $data
